$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'319.59"
$ws.Range("E2").Value = "'3.35%"
$ws.Range("D3").Value = "'41.44"
$ws.Range("E3").Value = "'1.40%"
$ws.Range("D4").Value = "'5.259"
$ws.Range("E4").Value = "'2.58%"
$ws.Range("D5").Value = "'0.07744"
$ws.Range("E5").Value = "'1.58%"
$ws.Range("D6").Value = "'1.754"
$ws.Range("E6").Value = "'9.22%"
$ws.Range("D7").Value = "'0.9461"
$ws.Range("E7").Value = "'4.21%"
$ws.Range("D9").Value = "'0.1236"
$ws.Range("E9").Value = "'-2.51%"
$ws.Range("D10").Value = "'0.1868"
$ws.Range("E10").Value = "'3.42%"
$ws.Range("D11").Value = "'0.09183"
$ws.Range("E11").Value = "'1.31%"
$ws.Range("D12").Value = "'0.04105"
$ws.Range("E12").Value = "'-5.21%"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'0.58%"
$ws.Range("D14").Value = "'0.001279"
$ws.Range("E14").Value = "'2.05%"
$ws.Range("D15").Value = "'0.005804"
$ws.Range("E15").Value = "'2.16%"
$ws.Range("D17").Value = "'3.351"
$ws.Range("E17").Value = "'-0.04%"
$ws.Range("D18").Value = "'4.355"
$ws.Range("E18").Value = "'1.48%"
$ws.Range("D20").Value = "'8.732"
$ws.Range("E20").Value = "'26.27%"
$ws.Range("E21").Value = "'-2.33%"
$ws.Range("D22").Value = "'0.2822"
$ws.Range("E22").Value = "'3.02%"
$ws.Range("D23").Value = "'0.04019"
$ws.Range("E23").Value = "'-0.66%"
$ws.Range("E24").Value = "'-0.19%"
$ws.Range("D25").Value = "'0.004121"
$ws.Range("E25").Value = "'1.81%"
$ws.Range("D26").Value = "'0.0001271"
$ws.Range("E26").Value = "'-0.11%"
$ws.Range("D38").Value = "'0.02566"
$ws.Range("E38").Value = "'6.20%"
$ws.Range("D39").Value = "'0.05338"
$ws.Range("E39").Value = "'1.98%"
$ws.Range("D40").Value = "'0.007778"
$ws.Range("E40").Value = "'-0.81%"
$ws.Range("E41").Value = "'1.24%"
$ws.Range("D42").Value = "'0.007043"
$ws.Range("E42").Value = "'3.57%"
$ws.Range("E43").Value = "'6.91%"
$ws.Range("D44").Value = "'0.008231"
$ws.Range("E44").Value = "'10.76%"
$ws.Range("D45").Value = "'0.3173"
$ws.Range("D46").Value = "'0.00006701"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("D48").Value = "'0.1975"
$ws.Range("E48").Value = "'23.96%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.11%"
